$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.566.79'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +5.30%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.371.99'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +10.83%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '257.89'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +11.10%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '624.22'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.67%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.22'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +12.98%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.387'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +4.48%  '

# Row 9
$ws.Range("E9").Value = '  +0.00%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.368.10'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +10.70%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.807'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.22%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.200'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.08%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '98.259.01'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +5.18%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.04'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +8.88%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000247'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.44%  '

# Row 16
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.995.98'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +10.71%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.51'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.26%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.371.08'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +10.06%  '

# Row 19
$ws.Range("E19").Value = '  +3.29%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.10'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +6.72%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '485.01'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +12.13%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.86'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.46%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000208'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +11.76%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.22'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +6.09%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.68'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +4.81%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.46'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.02%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.04'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +4.16%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.556.23'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +10.44%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.255'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.55%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.186'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +6.76%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.126'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +4.42%  '

# Row 33
$ws.Range("E33").Value = '  -9.21%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.33'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.58%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.31'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +9.32%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.44'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.36%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '518.76'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +14.02%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.151'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.96%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.95'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.36%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.91'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.06%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.447'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.54%  '

# Row 42
$ws.Range("E42").Value = '  +3.29%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.61'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.47%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.27'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +6.75%  '

# Row 45
$ws.Range("E45").Value = '  +0.01%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.776'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +18.07%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '160.85'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.23%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.92'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +6.94%  '

# Row 49
$ws.Range("E49").Value = '  +9.12%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.57'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.33%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.52'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +8.36%  '

Write-Host "Update complete"